# Revert "GUION 3 GRADO 8 DEFINITIVOS"
# This reverts the tracked-change deletions of trailing periods and
# restores the merged/cleaned-up "palabras clave" run, moving the
# _GoBack bookmark back to the top of the document.

$d = $word.ActiveDocument

# Work without generating new tracked changes.
$d.TrackRevisions = $false

# 1) The document contains 9 tracked deletions (all deleting a lone "."
#    character). The target revision keeps that text as normal,
#    un-tracked content, i.e. the deletions must be rejected (undone),
#    not accepted.
if ($d.Revisions.Count -gt 0) {
    $d.RejectAllRevisions()
}

# 2) Clean up / normalize the "Palabras clave" keyword list so the
#    comma-separated values each get a following space, and the stray
#    proofing-error markers collapse away as Word re-merges the runs.
$r = $d.Content
$r.Find.Execute(
    " rebelión,Túpac Amaru,Túpac Katari,comuneros,Manuela Beltrán,José Antonio Galán",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " rebelión, Túpac Amaru, Túpac Katari, comuneros, Manuela Beltrán, José Antonio Galán",
    2
) | Out-Null

# 2b) The 3 "Elige la respuesta correcta" sentences had their closing
#     period restored as a separate run by RejectAllRevisions above;
#     re-typing the full sentence merges it back into a single run,
#     matching how the original edit collapsed it.
$r2 = $d.Content
$r2.Find.Execute(
    "Luego de leer el texto, responde a la pregunta. Elige la respuesta correcta.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Luego de leer el texto, responde a la pregunta. Elige la respuesta correcta.",
    2
) | Out-Null

# 3) Move the _GoBack bookmark back to the empty paragraph right after
#    the title ("Ejercicio Genérico M5D: Test - con texto largo"),
#    which is the document's 2nd paragraph. Adding a new _GoBack
#    bookmark automatically removes whichever one previously existed
#    (further down the document).
$target = $d.Paragraphs.Item(2).Range
$d.Bookmarks.Add("_GoBack", $target) | Out-Null

# Restore the document's original "track changes" setting (it was on
# before this edit and the diff does not touch word/settings.xml).
$d.TrackRevisions = $true
